$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "97÷9=10, 7"; New = "15÷8=1, 7" },
    @{ Old = "19÷6=3, 1"; New = "77÷9=8, 5" },
    @{ Old = "46÷9=5, 1"; New = "97÷8=12, 1" },
    @{ Old = "35÷5=7, 0"; New = "73÷9=8, 1" },
    @{ Old = "75÷9=8, 3"; New = "30÷5=6, 0" },
    @{ Old = "61÷5=12, 1"; New = "86÷8=10, 6" },
    @{ Old = "25÷3=8, 1"; New = "44÷8=5, 4" },
    @{ Old = "75÷8=9, 3"; New = "96÷2=48, 0" },
    @{ Old = "91÷6=15, 1"; New = "20÷3=6, 2" },
    @{ Old = "23÷2=11, 1"; New = "54÷6=9, 0" },
    @{ Old = "82÷3=27, 1"; New = "66÷7=9, 3" },
    @{ Old = "68÷2=34, 0"; New = "95÷6=15, 5" },
    @{ Old = "96÷5=19, 1"; New = "64÷6=10, 4" },
    @{ Old = "36÷8=4, 4"; New = "49÷5=9, 4" },
    @{ Old = "45÷8=5, 5"; New = "13÷7=1, 6" },
    @{ Old = "64÷8=8, 0"; New = "10÷7=1, 3" },
    @{ Old = "49÷7=7, 0"; New = "95÷6=15, 5" },
    @{ Old = "21÷3=7, 0"; New = "87÷9=9, 6" },
    @{ Old = "64÷2=32, 0"; New = "13÷6=2, 1" },
    @{ Old = "13÷8=1, 5"; New = "87÷4=21, 3" },
    @{ Old = "43÷4=10, 3"; New = "59÷4=14, 3" },
    @{ Old = "50÷2=25, 0"; New = "32÷5=6, 2" },
    @{ Old = "87÷5=17, 2"; New = "98÷6=16, 2" },
    @{ Old = "94÷9=10, 4"; New = "38÷3=12, 2" },
    @{ Old = "77÷6=12, 5"; New = "34÷4=8, 2" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
